$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line")

# --- New "Notes" column header (I1), matching the existing header formatting ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Notes"

# --- New data rows 24-27: DE00-SE04 interconnection (TYNDP economic needs) ---
$ws.Range("A24").Value = "DE00-SE04"
$ws.Range("B24").Value = "Export Capacity"
$ws.Range("C24").Value = "Interconnection"
$ws.Range("D24").Value = "Distributed Energy"
$ws.Range("E24").Value = "Reference Grid"
$ws.Range("F24").Value = 2030
$ws.Range("G24").Value = 1984
$ws.Range("H24").Value = 1315
$ws.Range("I24").Value = "TYNDP economic needs"

$ws.Range("A25").Value = "DE00-SE04"
$ws.Range("B25").Value = "Import Capacity"
$ws.Range("C25").Value = "Interconnection"
$ws.Range("D25").Value = "Distributed Energy"
$ws.Range("E25").Value = "Reference Grid"
$ws.Range("F25").Value = 2030
$ws.Range("G25").Value = 1984
$ws.Range("H25").Value = -1315
$ws.Range("I25").Value = "TYNDP economic needs"

$ws.Range("A26").Value = "DE00-SE04"
$ws.Range("B26").Value = "Export Capacity"
$ws.Range("C26").Value = "Interconnection"
$ws.Range("D26").Value = "Distributed Energy"
$ws.Range("E26").Value = "Reference Grid"
$ws.Range("F26").Value = 2040
$ws.Range("G26").Value = 1984
$ws.Range("H26").Value = 1315
$ws.Range("I26").Value = "TYNDP economic needs"

$ws.Range("A27").Value = "DE00-SE04"
$ws.Range("B27").Value = "Import Capacity"
$ws.Range("C27").Value = "Interconnection"
$ws.Range("D27").Value = "Distributed Energy"
$ws.Range("E27").Value = "Reference Grid"
$ws.Range("F27").Value = 2040
$ws.Range("G27").Value = 1984
$ws.Range("H27").Value = -1315
$ws.Range("I27").Value = "TYNDP economic needs"

# --- New data rows 28-31: PL00-SE04 interconnection (TYNDP economic needs) ---
$ws.Range("A28").Value = "PL00-SE04"
$ws.Range("B28").Value = "Export Capacity"
$ws.Range("C28").Value = "Interconnection"
$ws.Range("D28").Value = "Distributed Energy"
$ws.Range("E28").Value = "Reference Grid"
$ws.Range("F28").Value = 2030
$ws.Range("G28").Value = 1984
$ws.Range("H28").Value = 1300
$ws.Range("I28").Value = "TYNDP economic needs"

$ws.Range("A29").Value = "PL00-SE04"
$ws.Range("B29").Value = "Import Capacity"
$ws.Range("C29").Value = "Interconnection"
$ws.Range("D29").Value = "Distributed Energy"
$ws.Range("E29").Value = "Reference Grid"
$ws.Range("F29").Value = 2030
$ws.Range("G29").Value = 1984
$ws.Range("H29").Value = -1300
$ws.Range("I29").Value = "TYNDP economic needs"

$ws.Range("A30").Value = "PL00-SE04"
$ws.Range("B30").Value = "Export Capacity"
$ws.Range("C30").Value = "Interconnection"
$ws.Range("D30").Value = "Distributed Energy"
$ws.Range("E30").Value = "Reference Grid"
$ws.Range("F30").Value = 2040
$ws.Range("G30").Value = 1984
$ws.Range("H30").Value = 1300
$ws.Range("I30").Value = "TYNDP economic needs"

$ws.Range("A31").Value = "PL00-SE04"
$ws.Range("B31").Value = "Import Capacity"
$ws.Range("C31").Value = "Interconnection"
$ws.Range("D31").Value = "Distributed Energy"
$ws.Range("E31").Value = "Reference Grid"
$ws.Range("F31").Value = 2040
$ws.Range("G31").Value = 1984
$ws.Range("H31").Value = -1300
$ws.Range("I31").Value = "TYNDP economic needs"

# --- Column I width (~24.74 chars in the original file's units) ---
$ws.Columns.Item(9).ColumnWidth = 23.8

# --- View: zoom to 100%, move the selection to I22 ---
$excel.ActiveWindow.Zoom = 100
$null = $ws.Range("I22").Select()
